$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 14

$ws.Cells.Item($row, 1).Value = 1
$ws.Cells.Item($row, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item($row, 3).Value = "Arica y Parinacota"

$ws.Cells.Item($row, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item($row, 4).Value = Get-Date -Year 2023 -Month 10 -Day 5 -Hour 0 -Minute 0 -Second 0

$ws.Cells.Item($row, 5).Value = 15
$ws.Cells.Item($row, 6).Value = 100112030
$ws.Cells.Item($row, 7).Value = "Poroto granado"
$ws.Cells.Item($row, 8).Value = "Sin especificar"
$ws.Cells.Item($row, 9).Value = "Primera"
$ws.Cells.Item($row, 10).Value = 1200
$ws.Cells.Item($row, 11).Value = 1600
$ws.Cells.Item($row, 12).Value = 1700
$ws.Cells.Item($row, 13).Value = 1650
$ws.Cells.Item($row, 14).Value = "`$/kilo"
$ws.Cells.Item($row, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item($row, 16).Value = 1650
$ws.Cells.Item($row, 17).Value = 1
$ws.Cells.Item($row, 18).Value = "Hortaliza"
